# Update the "Förändrad" (changed) date in column C for the existing data
# rows (2-43) from 45744 (2025-03-28) to 45745 (2025-03-29).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C43").Value2 = 45745

# The two newest report rows (44 "A 10579-2025" and 45 "A 13183-2025") were
# removed from the log, so delete those entire rows. Delete bottom-up so the
# remaining row numbers don't shift out from under us.
$ws.Rows.Item(45).Delete()
$ws.Rows.Item(44).Delete()

# Row 43 is now the last row; it should use the default (auto) row height
# instead of the explicit custom height it had before.
$ws.Rows.Item(43).AutoFit()
